$wb = $excel.ActiveWorkbook
$wsS0 = $wb.Worksheets.Item("scenarios")

# Re-apply the AutoFilter over the original data extent (A1:N49) BEFORE
# any new rows are written below it, otherwise the filter range silently
# grows to track the new bottom of the data.
$wsS0.Range("A1:N49").AutoFilter()

$names = $wb.Names
for ($i = 1; $i -le $names.Count(); $i++) {
  $n = $names.Item($i)
  if ($n.Name() -eq "scenarios!_FilterDatabase") {
    $n.RefersTo = "=scenarios!`$A`$1:`$N`$49"
  }
}

# ----------------------------------------------------------------------
# 1) descriptions sheet: re-purpose rows 15-20 and extend through row 31
#    with the new E-misspecification (M##) codes, then F/L rows that used
#    to live in 15-20 move down to 26-31.
# ----------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("descriptions")

$descRows = @(
  @(15, "E", 3, "internalCV"),
  @(16, "E", 10, "M50"),
  @(17, "E", 11, "M60"),
  @(18, "E", 12, "M70"),
  @(19, "E", 13, "M80"),
  @(20, "E", 14, "M90"),
  @(21, "E", 15, "M110"),
  @(22, "E", 16, "M120"),
  @(23, "E", 17, "M130"),
  @(24, "E", 18, "M140"),
  @(25, "E", 19, "M150"),
  @(26, "F", 0, "constant"),
  @(27, "F", 1, "contrast"),
  @(28, "F", 2, "increase"),
  @(29, "L", 10, "fish length"),
  @(30, "L", 30, "fish & surv length"),
  @(31, "L", 31, "fish & .5surv length")
)

foreach ($row in $descRows) {
  $r = $row[0]
  $wsD.Cells.Item($r, 1).Value = $row[1]
  $wsD.Cells.Item($r, 2).Value = $row[2]
  $wsD.Cells.Item($r, 4).Value = $row[3]
  $wsD.Cells.Item($r, 3).Formula = "=CONCATENATE(A" + $r + ",B" + $r + ")"
}

# ----------------------------------------------------------------------
# 2) scenarios sheet: append 20 new scenario rows (50-69) covering the
#    new E10..E19 misspecification cases.
# ----------------------------------------------------------------------
$wsS = $wb.Worksheets.Item("scenarios")

$scenRows = @(
  @(50, "A10", "C0", "D0", "L10", "E10", "F1", "I0"),
  @(51, "A10", "C0", "D10", "L10", "E10", "F1", "I0"),
  @(52, "A10", "C0", "D0", "L10", "E11", "F1", "I0"),
  @(53, "A10", "C0", "D10", "L10", "E11", "F1", "I0"),
  @(54, "A10", "C0", "D0", "L10", "E12", "F1", "I0"),
  @(55, "A10", "C0", "D10", "L10", "E12", "F1", "I0"),
  @(56, "A10", "C0", "D0", "L10", "E13", "F1", "I0"),
  @(57, "A10", "C0", "D10", "L10", "E13", "F1", "I0"),
  @(58, "A10", "C0", "D0", "L10", "E14", "F1", "I0"),
  @(59, "A10", "C0", "D10", "L10", "E14", "F1", "I0"),
  @(60, "A10", "C0", "D0", "L10", "E15", "F1", "I0"),
  @(61, "A10", "C0", "D10", "L10", "E15", "F1", "I0"),
  @(62, "A10", "C0", "D0", "L10", "E16", "F1", "I0"),
  @(63, "A10", "C0", "D10", "L10", "E16", "F1", "I0"),
  @(64, "A10", "C0", "D0", "L10", "E17", "F1", "I0"),
  @(65, "A10", "C0", "D10", "L10", "E17", "F1", "I0"),
  @(66, "A10", "C0", "D0", "L10", "E18", "F1", "I0"),
  @(67, "A10", "C0", "D10", "L10", "E18", "F1", "I0"),
  @(68, "A10", "C0", "D0", "L10", "E19", "F1", "I0"),
  @(69, "A10", "C0", "D10", "L10", "E19", "F1", "I0")
)

foreach ($row in $scenRows) {
  $r = $row[0]
  $wsS.Cells.Item($r, 1).Value = $row[1]
  $wsS.Cells.Item($r, 2).Value = $row[2]
  $wsS.Cells.Item($r, 3).Value = $row[3]
  $wsS.Cells.Item($r, 4).Value = $row[4]
  $wsS.Cells.Item($r, 5).Value = $row[5]
  $wsS.Cells.Item($r, 6).Value = $row[6]
  $wsS.Cells.Item($r, 7).Value = $row[7]

  $wsS.Cells.Item($r, 9).Formula = "=LOOKUP(A" + $r + ",descriptions!`$C:`$C,descriptions!`$D:`$D)"
  $wsS.Cells.Item($r, 10).Formula = "=LOOKUP(B" + $r + ",descriptions!`$C:`$C,descriptions!`$D:`$D)"
  $wsS.Cells.Item($r, 11).Formula = "=LOOKUP(C" + $r + ",descriptions!`$C:`$C,descriptions!`$D:`$D)"
  $wsS.Cells.Item($r, 12).Formula = "=LOOKUP(D" + $r + ",descriptions!`$C:`$C,descriptions!`$D:`$D)"
  $wsS.Cells.Item($r, 13).Formula = "=LOOKUP(E" + $r + ",descriptions!`$C:`$C,descriptions!`$D:`$D)"
  $wsS.Cells.Item($r, 14).Formula = "=LOOKUP(F" + $r + ",descriptions!`$C:`$C,descriptions!`$D:`$D)"
}

# H column: rows 50-53 get their own formula (matches how the original
# author typed/filled them one at a time), rows 54-69 were filled in one
# drag so they share a single formula definition.
$wsS.Range("H50").Formula = "=CONCATENATE(A50,""-"",B50,""-"",C50,""-"",D50,""-"",E50,""-"",F50,""-"",G50,""-"")"
$wsS.Range("H51").Formula = "=CONCATENATE(A51,""-"",B51,""-"",C51,""-"",D51,""-"",E51,""-"",F51,""-"",G51,""-"")"
$wsS.Range("H52").Formula = "=CONCATENATE(A52,""-"",B52,""-"",C52,""-"",D52,""-"",E52,""-"",F52,""-"",G52,""-"")"
$wsS.Range("H53").Formula = "=CONCATENATE(A53,""-"",B53,""-"",C53,""-"",D53,""-"",E53,""-"",F53,""-"",G53,""-"")"
$wsS.Range("H54:H69").Formula = "=CONCATENATE(A54,""-"",B54,""-"",C54,""-"",D54,""-"",E54,""-"",F54,""-"",G54,""-"")"

# ----------------------------------------------------------------------
# 3) Restore the on-screen selections recorded in the saved view state.
# ----------------------------------------------------------------------
$wsD.Rows.Item(26).Select()

$wsS.Activate()
$wsS.Range("N63").Select()
